# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
#
# The sheet already has header + stat columns through AC ("Unnamed: 28").
# We append three new columns - AD "Wins", AE "Losses", AF "Ties" - with a
# header style matching the existing headers (bold, centered, bordered),
# and fill every player row (2-51) with the team's season record:
# 66 wins, 96 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns AD, AE, AF.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the other header cells (bold font + border) by copying
# the formatting from the neighboring header cell (AC1) onto the new ones.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record for every data row.
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 66
    $ws.Cells.Item($r, 31).Value = 96
    $ws.Cells.Item($r, 32).Value = 0
}
